$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 2.341355666666667
$ws.Range("N2").Value = 7.024067000000001
$ws.Range("O2").Value = 0.03973512964576821
$ws.Range("P2").Value = 0.0397351296457682
$ws.Range("Q2").Value = 0.491761954737
$ws.Range("R2").Value = 4.425857592633
$ws.Range("S2").Value = 0.01239343341048897
$ws.Range("T2").Value = 0.01239343341048897

# Row 3
$ws.Range("O3").Value = 0.5779093692199981
$ws.Range("P3").Value = 0.5779093692199981
$ws.Range("Q3").Value = 7.152206211531
$ws.Range("R3").Value = 64.369855903779
$ws.Range("S3").Value = 0.1802506081791158
$ws.Range("T3").Value = 0.1802506081791158

# Row 4
$ws.Range("O4").Value = 0.3823555011342337
$ws.Range("P4").Value = 0.3823555011342337
$ws.Range("S4").Value = 0.1192571280045122
$ws.Range("T4").Value = 0.1192571280045122

# Row 5
$ws.Range("I5").Value = 0.688098830405883
$ws.Range("J5").Value = 0.688098830405883
$ws.Range("M5").Value = 2.341355666666667
$ws.Range("N5").Value = 7.024067000000001
$ws.Range("O5").Value = 0.03973512964576821
$ws.Range("P5").Value = 0.0397351296457682
$ws.Range("Q5").Value = 1.084897585773667
$ws.Range("R5").Value = 9.764078271963001
$ws.Range("S5").Value = 0.02734169623527924
$ws.Range("T5").Value = 0.02734169623527923

# Row 6
$ws.Range("I6").Value = 0.688098830405883
$ws.Range("J6").Value = 0.688098830405883
$ws.Range("O6").Value = 0.5779093692199981
$ws.Range("P6").Value = 0.5779093692199981
$ws.Range("S6").Value = 0.3976587610408823
$ws.Range("T6").Value = 0.3976587610408823

# Row 7
$ws.Range("I7").Value = 0.688098830405883
$ws.Range("J7").Value = 0.688098830405883
$ws.Range("O7").Value = 0.3823555011342337
$ws.Range("P7").Value = 0.3823555011342337
$ws.Range("R7").Value = 93.95587919487002
$ws.Range("S7").Value = 0.2630983731297215
$ws.Range("T7").Value = 0.2630983731297215
